$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch")
# ("ch" is also the workbook's one and only / active sheet: $wb.ActiveSheet)

# Update the "Confirm to remove ..." translation text (the/this)
$ws.Range("C2").Value = "Confirm to remove this security group"

# Widen column C to fit the longer text (target stored width ~39.21875 chars;
# COM ColumnWidth is quantised to 1/7 character units by the engine, so land
# on the closest representable value)
$ws.Columns.Item(3).ColumnWidth = 38.5714285714286

# Set up the page for printing (adds <pageSetup .../> to the sheet)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the active selection to C16
$ws.Range("C16").Select() | Out-Null
